$wb = $excel.ActiveWorkbook

# Add the new "Summary" sheet right after the existing data sheet.
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $firstSheet)
$ws.Name = "Summary"

# Header row.
$ws.Cells.Item(1, 1).Value = "Domain"
$ws.Cells.Item(1, 2).Value = "Total Found"

# Stack-ranked counts of each target account domain found in the
# DataMiner search results (descending by count).
$summaryData = @(
    @("gdit.com", 135),
    @("ngc.com", 58),
    @("homedepot.com", 44),
    @("ihg.com", 37),
    @("fanniemae.com", 35),
    @("delta.com", 35),
    @("lfg.com", 35),
    @("mckesson.com", 29),
    @("raymondjames.com", 29),
    @("freddiemac.com", 22),
    @("rccl.com", 19),
    @("navyfederal.org", 19),
    @("manh.com", 17),
    @("macys.com", 16),
    @("fnf.com", 13),
    @("autotrader.com", 12),
    @("na.ko.com", 10),
    @("blackbaud.com", 10),
    @("advance-auto.com", 10),
    @("sungardas.com", 9),
    @("publix.com", 9),
    @("lowes.com", 9),
    @("fiserv.com", 7),
    @("ncr.com", 7),
    @("harris.com", 7),
    @("email.chop.edu", 6),
    @("marriott.com", 6),
    @("bcbsnc.com", 6),
    @("geico.com", 6),
    @("bcbsfl.com", 5),
    @("officedepot.com", 5),
    @("equifax.com", 5),
    @("ugcorp.com", 4),
    @("vanguard.com", 4),
    @("catalinamarketing.com", 4),
    @("fpl.com", 4),
    @("sas.com", 4),
    @("jmfamily.com", 3),
    @("ge.com", 3),
    @("duke-energy.com", 3),
    @("dollartree.com", 3),
    @("dell.com", 3),
    @("imshealth.com", 3),
    @("neustar.biz", 3),
    @("labcorp.com", 3),
    @("carefirst.com", 3),
    @("aarp.org", 2),
    @("citrix.com", 2),
    @("dominionenterprises.com", 2),
    @("underarmour.com", 2),
    @("effem.com", 2),
    @("hilton.com", 2),
    @("moffitt.org", 2),
    @("hcsc.net", 2),
    @("amerisourcebergen.com", 2),
    @("bdpint.com", 1),
    @("comscore.com", 1),
    @("avidxchange.com", 1),
    @("csx.com", 1),
    @("travelport.com", 1),
    @("sita.aero", 1),
    @("blackboard.com", 1),
    @("verisign.com", 1),
    @("wfu.edu", 1),
    @("hersheys.com", 1),
    @("healthesystems.com", 1),
    @("synchronoss.com", 1),
    @("ups.com", 0),
    @("danaher.com", 0),
    @("chicos.com", 0),
    @("southernco.com", 0),
    @("markelcorp.com", 0),
    @("genworth.com", 0),
    @("inovalon.com", 0),
    @("sbgnet.com", 0),
    @("hanloninvest.com", 0),
    @("merck.com", 0),
    @("nielsen.com", 0),
    @("altisource.com", 0),
    @("microstrategy.com", 0),
    @("freedommortgage.com", 0),
    @("masonite.com", 0),
    @("autonation.com", 0),
    @("astrazeneca.com", 0),
    @("sykes.com", 0),
    @("slhn.org", 0),
    @("subaru.com", 0),
    @("footballfanatics.com", 0),
    @("carmax.com", 0),
    @("bbandt.com", 0),
    @("aflac.com", 0),
    @("baycare.org", 0),
    @("troweprice.com", 0),
    @("vertexinc.com", 0),
    @("iassoftware.com", 0),
    @("tsys.com", 0),
    @("carnival.com", 0),
    @("darden.com", 0),
    @("ahss.org", 0),
    @("syniverse.com", 0),
    @("nascar.com", 0),
    @("wellcare.com", 0),
    @("mohawkind.com", 0),
    @("transcore.com", 0),
    @("carecorenational.com", 0),
    @("usa.dupont.com", 0),
    @("ultimatesoftware.com", 0),
    @("fticonsulting.com", 0),
    @("benefitfocus.com", 0),
    @("praintl.com", 0)
)

for ($i = 0; $i -lt $summaryData.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $summaryData[$i][0]
    $ws.Cells.Item($r, 2).Value = $summaryData[$i][1]
}

# Match the column sizing used in the source report (best-fit to content).
$ws.Columns.Item(1).ColumnWidth = 21.065104166666668
$ws.Columns.Item(2).ColumnWidth = 9.803385416666666

$wb.Worksheets.Item(1).Activate()
